# Append 21 new coded-segment rows (rows 219-239) to Sheet1, reproducing the
# "Run through with latest mex files" commit: new rows tagged with coder
# "chen" and timestamps from 1/31/19, reusing existing categorical text
# (document/code/segment values) that is already present elsewhere in the
# sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: D, F, G, I, J, K, M  (A/B/C/E/H/L are constant across all new rows)
$rowsData = @(
    @("6314", "4: 3444", "4: 3447", "_x0002_256", "4", "9.1907541013740175E-3", "1/31/19 13:11:03"),
    @("6314", "4: 3527", "4: 3530", "_x0002_256", "4", "9.1907541013740175E-3", "1/31/19 13:11:08"),
    @("6314", "4: 3574", "4: 3577", "_x0002_256", "4", "9.1907541013740175E-3", "1/31/19 13:11:13"),
    @("6314", "4: 3623", "4: 3626", "_x0002_256", "4", "9.1907541013740175E-3", "1/31/19 13:11:19"),
    @("6314", "4: 3600", "4: 3603", "_x0002_256", "4", "9.1907541013740175E-3", "1/31/19 13:11:24"),
    @("6314", "4: 3675", "4: 3677", "_x0002_32", "3", "6.8930655760305135E-3", "1/31/19 13:11:44"),
    @("6314", "4: 3248", "4: 3251", "_x0002_256", "4", "9.1907541013740175E-3", "1/31/19 13:12:01"),
    @("6314", "4: 3266", "4: 3269", "_x0002_256", "4", "9.1907541013740175E-3", "1/31/19 13:12:04"),
    @("6314", "4: 3290", "4: 3293", "_x0002_256", "4", "9.1907541013740175E-3", "1/31/19 13:12:14"),
    @("11395", "3: 3624", "3: 3626", "512", "3", "1.6090104585679808E-2", "1/31/19 13:33:01"),
    @("11395", "3: 3648", "3: 3651", ">512", "4", "2.145347278090641E-2", "1/31/19 13:33:05"),
    @("11395", "3: 3695", "3: 3698", ">512", "4", "2.145347278090641E-2", "1/31/19 13:33:10"),
    @("11395", "3: 3718", "3: 3720", "256", "3", "1.6090104585679808E-2", "1/31/19 13:33:25"),
    @("11395", "3: 3509", "3: 3512", ">512", "4", "2.145347278090641E-2", "1/31/19 13:33:40"),
    @("11395", "3: 3534", "3: 3537", ">512", "4", "2.145347278090641E-2", "1/31/19 13:33:44"),
    @("11395", "3: 3595", "3: 3598", ">512", "4", "2.145347278090641E-2", "1/31/19 13:34:03"),
    @("13370", "2: 5938", "2: 5940", "_x0004_32", "3", "8.0394468860542391E-3", "1/31/19 13:35:10"),
    @("13370", "2: 5945", "2: 5946", "_x0004_8", "2", "5.359631257369493E-3", "1/31/19 13:35:14"),
    @("13370", "3: 5704", "3: 5705", "_x0004_8", "2", "5.359631257369493E-3", "1/31/19 13:35:44"),
    @("13370", "2: 5924", "2: 5926", "_x0004_16", "3", "8.0394468860542391E-3", "1/31/19 13:35:56"),
    @("13370", "2: 5928", "2: 5930", "_x0004_16", "3", "8.0394468860542391E-3", "1/31/19 13:36:01")
)

$startRow = 219

for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $r = $startRow + $i
    $d = $rowsData[$i]

    # Clone the full row formatting (fills/borders/fonts/alignment/number
    # formats) from row 2, which carries the same per-column style pattern
    # (A=6,B=2,C=1,D=2,E=1,F=1,G=1,H=3,I=2,J=3,K=4,L=1,M=1) as every data row.
    $src = $ws.Range("A2:M2")
    $dst = $ws.Range("A" + $r + ":M" + $r)
    $src.Copy($dst)
    $ws.Rows.Item($r).RowHeight = 16

    # Constant columns for every appended row.
    $ws.Cells.Item($r, 1).Value = "$([char]9679)"   # A: bullet marker (same as every other row)
    $ws.Cells.Item($r, 2).Value = ""                 # B
    $ws.Cells.Item($r, 3).Value = ""                 # C
    $ws.Cells.Item($r, 5).Value = "MIC"               # E

    # Text columns that might otherwise be auto-converted to numbers/dates by
    # Excel on assignment - force Text format first so they round-trip as
    # shared strings exactly like the source data.
    $dCell = $ws.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $d[0]

    $ws.Cells.Item($r, 6).Value = $d[1]   # F
    $ws.Cells.Item($r, 7).Value = $d[2]   # G

    $ws.Cells.Item($r, 8).Value = 0.0     # H (numeric)

    $iCell = $ws.Cells.Item($r, 9)
    $iCell.NumberFormat = "@"
    $iCell.Value = $d[3]

    $ws.Cells.Item($r, 10).Value = [double]$d[4]   # J (numeric count)
    $ws.Cells.Item($r, 11).Value = [double]$d[5]   # K (numeric fraction)

    $lCell = $ws.Cells.Item($r, 12)
    $lCell.NumberFormat = "@"
    $lCell.Value = "chen"                 # L: coder name

    $mCell = $ws.Cells.Item($r, 13)
    $mCell.NumberFormat = "@"
    $mCell.Value = $d[6]                  # M: timestamp
}

Write-Host "Appended" $rowsData.Count "rows starting at" $startRow
